$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "notes" queue in column J shifts up: the two oldest entries (J9, J12 -
# along with already-empty J9's sibling) are consumed/completed, matching the
# commit "CarSearchController now uses stored procedure - removed old search
# utility". Every remaining J note slides up two rows, and a freshly queued
# note ("Add Exception handling - add to db", formerly a stray item at H52)
# is appended at the bottom of the queue.

# Capture current (pre-edit) J-column note text before we start overwriting,
# so the "shift" reads cleanly from a stable snapshot.
$j16 = $ws.Range("J16").Value()
$j18 = $ws.Range("J18").Value()
$j20 = $ws.Range("J20").Value()
$j21 = $ws.Range("J21").Value()
$j22 = $ws.Range("J22").Value()
$h52 = $ws.Range("H52").Value()

# Row 9: oldest note is fully done -> remove (Clear, not ClearContents, so the
# cell element disappears entirely rather than lingering with its old style).
$ws.Range("J9").Clear()

# Row 10: now holds what used to be two rows below's content position-wise,
# i.e. the note that used to sit at J12.
$ws.Range("J10").Value() = $ws.Range("J12").Value()

# Row 12: cleared (its content moved up into J10).
$ws.Range("J12").Clear()

# Row 14: gains a new note cell (style matches the other queued notes),
# carrying what used to be at J16.
$ws.Range("J14").Value() = $j16
$ws.Range("J14").WrapText = $true

# Row 16: now carries what used to be at J18.
$ws.Range("J16").Value() = $j18

# Row 18: now carries what used to be at J20.
$ws.Range("J18").Value() = $j20

# Row 19: gains a new note cell, carrying what used to be at J21.
$ws.Range("J19").Value() = $j21
$ws.Range("J19").WrapText = $true

# Row 20: now carries what used to be at J22.
$ws.Range("J20").Value() = $j22

# Row 21 & 22: cleared (their content moved up) - Clear() so the now-unused
# note style doesn't leave a bare styled cell behind.
$ws.Range("J21").Clear()
$ws.Range("J22").Clear()

# Row 23: newly queued note, moved here from the orphan H52 entry below.
$ws.Range("J23").Value() = $h52

# The item that used to live at H52 (standalone backlog entry) is now tracked
# in the J queue instead, so remove it from the bottom mini-table - this also
# drops the row from the sheet's used range.
$ws.Range("H52").Clear()

# Reflect the new last-used selection.
$ws.Range("J28").Select()
